# Apply automatic electricity spot price update for row 2 of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45890

$ws.Range("B2").Value = 95.28
$ws.Range("C2").Value = 88.88
$ws.Range("D2").Value = 85.26000000000001
$ws.Range("E2").Value = 83.68000000000001
$ws.Range("F2").Value = 82.09999999999999
$ws.Range("G2").Value = 84
$ws.Range("H2").Value = 85.81999999999999
$ws.Range("I2").Value = 101.12
$ws.Range("J2").Value = 85.26000000000001
$ws.Range("K2").Value = 50.2
$ws.Range("L2").Value = 25.67
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 7.89
$ws.Range("O2").Value = 5.11
$ws.Range("P2").Value = 5.11
$ws.Range("Q2").Value = 4.31
$ws.Range("R2").Value = 5.79
$ws.Range("S2").Value = 6
$ws.Range("T2").Value = 25
$ws.Range("U2").Value = 50
$ws.Range("V2").Value = 77.84
$ws.Range("W2").Value = 83.26000000000001
$ws.Range("X2").Value = 75
$ws.Range("Y2").Value = 71.59999999999999
$ws.Range("Z2").Value = 53.92

$ws.Range("AA2").Value = "0h-4h"
$ws.Range("AB2").Value = 88.28

$ws.Range("AD2").Value = 93.47

$ws.Range("AF2").Value = 92.08
$ws.Range("AG2").Value = "9h-19h"
